# Extent Report Feature:
# Split the old "invalidCountryandCourse" / "invalidCollegeMajorGPA" sheets
# into five focused single-invalid-field sheets:
#   invalidCountry, invalidCourse, invalidCollege, invalidMajor, invalidGPA

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122
$xlCenter = -4108

# ---------------------------------------------------------------------
# 1) Insert the brand-new "invalidCountry" sheet right after mastersDetails
# ---------------------------------------------------------------------
$masters = $wb.Worksheets.Item("mastersDetails")
$wsCountry = $wb.Worksheets.Add($null, $masters)
$wsCountry.Name = "invalidCountry"

$wsCountry.Range("A1").Value = "Country "
$masters.Range("A1").Copy()
$wsCountry.Range("A1").PasteSpecial($xlPasteFormats)

$wsCountry.Range("A2").Value = 653
$wsCountry.Range("A3").Value = '*$&'
$wsCountry.Range("A4").Value = ''
$wsCountry.Range("A2:A4").HorizontalAlignment = $xlCenter

$wsCountry.Range("D6").Select()

# ---------------------------------------------------------------------
# 2) Turn the old "invalidCountryandCourse" sheet into "invalidCourse"
#    (same sheet identity, rows 5-7 dropped, rows 2-4 rewritten)
# ---------------------------------------------------------------------
$wsCourse = $wb.Worksheets.Item("invalidCountryandCourse")
$wsCourse.Rows("5:7").Delete()
$wsCourse.Name = "invalidCourse"

$wsCourse = $wb.Worksheets.Item("invalidCourse")
$wsCourse.Range("A2").Value = "United Kingdom"
$wsCourse.Range("B2").Value = 39421

$wsCourse.Range("A3").Value = "United Kingdom"
$wsCourse.Range("B3").Value = '*$&'

$wsCourse.Range("A4").Value = "United Kingdom"
$wsCourse.Range("B4").Value = ''

$wsCourse.Range("A1:A2").Select()

# ---------------------------------------------------------------------
# 3) Insert the brand-new "invalidMajor" sheet right after invalidCourse
# ---------------------------------------------------------------------
$wsCourse = $wb.Worksheets.Item("invalidCourse")
$wsMajor = $wb.Worksheets.Add($null, $wsCourse)
$wsMajor.Name = "invalidMajor"

$wsMajor.Range("A1").Value = "Country "
$wsMajor.Range("B1").Value = "Course"
$wsMajor.Range("C1").Value = "College"
$wsMajor.Range("D1").Value = "Major"
$masters.Range("A1:D1").Copy()
$wsMajor.Range("A1:D1").PasteSpecial($xlPasteFormats)

$wsMajor.Range("A2").Value = "United Kingdom"
$wsMajor.Range("B2").Value = "Computer Science"
$wsMajor.Range("C2").Value = "Nirma University"

$wsMajor.Range("A3").Value = "United Kingdom"
$wsMajor.Range("B3").Value = "Computer Science"
$wsMajor.Range("C3").Value = "Nirma University"
$wsMajor.Range("D3").Value = '*^$('

$wsMajor.Range("A4").Value = "United Kingdom"
$wsMajor.Range("B4").Value = "Computer Science"
$wsMajor.Range("C4").Value = "Nirma University"
$wsMajor.Range("D4").Value = 9834

$wsMajor.Range("A4").Select()

# ---------------------------------------------------------------------
# 4) Insert the brand-new "invalidCollege" sheet between invalidCourse
#    and invalidMajor
# ---------------------------------------------------------------------
$wsCourse = $wb.Worksheets.Item("invalidCourse")
$wsCollege = $wb.Worksheets.Add($null, $wsCourse)
$wsCollege.Name = "invalidCollege"

$wsCollege.Range("A1").Value = "Country "
$wsCollege.Range("B1").Value = "Course"
$wsCollege.Range("C1").Value = "College"
$masters.Range("A1:C1").Copy()
$wsCollege.Range("A1:C1").PasteSpecial($xlPasteFormats)

$wsCollege.Range("A2").Value = "United Kingdom"
$wsCollege.Range("B2").Value = "Computer Science"

$wsCollege.Range("A3").Value = "United Kingdom"
$wsCollege.Range("B3").Value = "Computer Science"
$wsCollege.Range("C3").Value = 865

$wsCollege.Range("A4").Value = "United Kingdom"
$wsCollege.Range("B4").Value = "Computer Science"
$wsCollege.Range("C4").Value = '&$%*'

$wsCollege.PageSetup.Orientation = 1
$wsCollege.Range("C4").Select()

# ---------------------------------------------------------------------
# 5) Turn the old "invalidCollegeMajorGPA" sheet into "invalidGPA"
#    (same sheet identity, rows 6-10 dropped, rows 2-5 rewritten,
#     columns extended from College/Major/GPA to the full 5-column set)
# ---------------------------------------------------------------------
$wsGPA = $wb.Worksheets.Item("invalidCollegeMajorGPA")
$wsGPA.Rows("6:10").Delete()
$wsGPA.Columns("C:C").Delete()
$wsGPA.Name = "invalidGPA"

$wsGPA = $wb.Worksheets.Item("invalidGPA")
$wsGPA.Range("A1").Value = "Country "
$wsGPA.Range("B1").Value = "Course"
$wsGPA.Range("C1").Value = "College"
$wsGPA.Range("D1").Value = "Major"
$wsGPA.Range("E1").Value = "GPA"
$masters.Range("A1:E1").Copy()
$wsGPA.Range("A1:E1").PasteSpecial($xlPasteFormats)

$wsGPA.Range("A2").Value = "United Kingdom"
$wsGPA.Range("B2").Value = "Computer Science"
$wsGPA.Range("C2").Value = "Nirma University"
$wsGPA.Range("D2").Value = "Electronics and Communication"
$wsGPA.Range("E2").Value = 'iwm'

$wsGPA.Range("A3").Value = "United Kingdom"
$wsGPA.Range("B3").Value = "Computer Science"
$wsGPA.Range("C3").Value = "Nirma University"
$wsGPA.Range("D3").Value = "Electronics and Communication"
$wsGPA.Range("E3").Value = '*^%('

$wsGPA.Range("A4").Value = "United Kingdom"
$wsGPA.Range("B4").Value = "Computer Science"
$wsGPA.Range("C4").Value = "Nirma University"
$wsGPA.Range("D4").Value = "Electronics and Communication"
$wsGPA.Range("E4").Value = ''

$wsGPA.Range("A5").Value = "United Kingdom"
$wsGPA.Range("B5").Value = "Computer Science"
$wsGPA.Range("C5").Value = "Nirma University"
$wsGPA.Range("D5").Value = "Electronics and Communication"
$wsGPA.Range("E5").Value = 23

$wsGPA.Range("D16").Select()

# ---------------------------------------------------------------------
# 6) Update mastersDetails selection and finish on invalidMajor as the
#    active tab
# ---------------------------------------------------------------------
$masters = $wb.Worksheets.Item("mastersDetails")
$masters.Activate()
$masters.Range("B12").Select()

$wb.Worksheets.Item("invalidMajor").Activate()
